$d = $word.ActiveDocument

# 1. "Raluca, " -> "Raluca-Anamaria, "
$d.Content.Find.Execute("Raluca, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Raluca-Anamaria, ", 2)

# 2. "portofolio" -> "portfolio" (two occurrences)
$d.Content.Find.Execute("portofolio", $true, $false, $false, $false, $false,
                         $true, 1, $false, "portfolio", 2)

Write-Output "done"
